# T1649_AdditionalClientsSubjectRequired.xlsx
# "updates after SF Referesh(Except EventExpenase)"
#
# Changes applied (AddOpportunity sheet unless noted):
#  1. Shared-string content "Tec Alliance" -> "CB Alliance" — cells O2 and
#     AG2 (both "Search" columns) get the new client name.
#  2. Clean up redundant cell formatting left over from earlier edits:
#       - R1 / S1 (bold header cells that also carried a stray "apply
#         number format" flag) are normalised to the plain bold style
#         already used by the rest of row 1 (same style as A1).
#       - A2 / B2 / R2 / S2 (data cells that also carried a stray "apply
#         number format" flag) are normalised to the default/general
#         style already used by the rest of row 2 (same style as C2).
#  3. Add two new (wider) helper columns R and S so their text fits,
#     and move the sheet's scroll position / active selection further
#     to the right (to around column AK / cell BM2) to reflect the
#     newly added columns further out in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddOpportunity")

# 1. Update the "Tec Alliance" -> "CB Alliance" text (both copies).
$ws.Range("O2").Value = "CB Alliance"
$ws.Range("AG2").Value = "CB Alliance"

# 2. Normalise formatting: R1/S1 should match the bold style already used
#    by A1 (and the rest of the header row); A2/B2/R2/S2 should match the
#    plain/default style already used by C2 (and the rest of row 2).
$ws.Range("A1").Copy() | Out-Null
$ws.Range("R1:S1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("C2").Copy() | Out-Null
$ws.Range("A2:B2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("R2:S2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

# 3. Widen the R/S columns and update the view's selection/scroll position.
$ws.Columns.Item(18).ColumnWidth = 15.25
$ws.Columns.Item(19).ColumnWidth = 16.59

$ws.Range("BM2").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 37
